$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.210126638412476
$ws.Range("B1").Value = 2.437742233276367
$ws.Range("C1").Value = 4.689794063568115
$ws.Range("D1").Value = 2.595423936843872
$ws.Range("E1").Value = 1.089164614677429
